$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F (想去人数 / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 273
$ws1.Range("F3").Value = 160
$ws1.Range("F4").Value = 1964
$ws1.Range("F5").Value = 1608
$ws1.Range("F6").Value = 288
$ws1.Range("F7").Value = 75
$ws1.Range("F8").Value = 605
$ws1.Range("F9").Value = 139

# Sheet "全部类型" (all types) - update column F (想去人数 / interested count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 273
$ws4.Range("F3").Value = 160
$ws4.Range("F4").Value = 1964
$ws4.Range("F5").Value = 1608
$ws4.Range("F6").Value = 288
$ws4.Range("F8").Value = 75
$ws4.Range("F9").Value = 605
$ws4.Range("F10").Value = 139
